$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC (sheet1.xml) ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 357.4
$ws.Range("I12").Value = 277.66666
$ws.Range("K12").Value = 277.66666
$ws.Range("M12").Value = -107.66666
$ws.Range("H17").Value = 178987.53
$ws.Range("J17").Value = 178987.53
$ws.Range("L17").Value = 536962.59
$ws.Range("N17").Value = -537298.59
$ws.Range("H28").Value = 1289.4286
$ws.Range("I28").Value = 782.75
$ws.Range("K28").Value = 782.75
$ws.Range("M28").Value = -297.75
$ws.Range("H32").Value = 10032.2
$ws.Range("I32").Value = 9687.125
$ws.Range("J32").Value = 10426.571
$ws.Range("K32").Value = 9687.125
$ws.Range("L32").Value = 10426.571
$ws.Range("M32").Value = -9361.125
$ws.Range("N32").Value = -11078.571
$ws.Range("H51").Value = 2930.4
$ws.Range("I51").Value = 1622.5
$ws.Range("J51").Value = 3406
$ws.Range("K51").Value = 1622.5
$ws.Range("L51").Value = 3406
$ws.Range("M51").Value = -1138.5
$ws.Range("N51").Value = -4374
$ws.Range("H62").Value = 115701.89
$ws.Range("I62").Value = 129288.75
$ws.Range("K62").Value = 129288.75
$ws.Range("M62").Value = -128664.75
$ws.Range("H64").Value = 5491.5
$ws.Range("J64").Value = 6738
$ws.Range("L64").Value = 6738
$ws.Range("N64").Value = -7234
$ws.Range("H65").Value = 115701.89
$ws.Range("I65").Value = 129288.75
$ws.Range("K65").Value = 646443.75
$ws.Range("M65").Value = -643323.75
$ws.Range("H67").Value = 5491.5
$ws.Range("J67").Value = 6738
$ws.Range("L67").Value = 6738
$ws.Range("N67").Value = -8454
$ws.Range("H70").Value = 84864.16
$ws.Range("I70").Value = 510000
$ws.Range("J70").Value = 7566.727
$ws.Range("K70").Value = 1530000
$ws.Range("L70").Value = 22700.181
$ws.Range("M70").Value = -1529730
$ws.Range("N70").Value = -23240.181
$ws.Range("H73").Value = 84864.16
$ws.Range("I73").Value = 510000
$ws.Range("J73").Value = 7566.727
$ws.Range("K73").Value = 1530000
$ws.Range("L73").Value = 22700.181
$ws.Range("M73").Value = -1529064
$ws.Range("N73").Value = -24572.181
$ws.Range("H80").Value = 2573.2942
$ws.Range("I80").Value = 2311.8333
$ws.Range("J80").Value = 2715.9092
$ws.Range("K80").Value = 6935.499899999999
$ws.Range("L80").Value = 8147.7276
$ws.Range("M80").Value = -5937.499899999999
$ws.Range("N80").Value = -10143.7276
$ws.Range("H83").Value = 2573.2942
$ws.Range("I83").Value = 2311.8333
$ws.Range("J83").Value = 2715.9092
$ws.Range("K83").Value = 20806.4997
$ws.Range("L83").Value = 24443.1828
$ws.Range("M83").Value = -15814.4997
$ws.Range("N83").Value = -34427.1828
$ws.Range("H98").Value = 2284.8845
$ws.Range("I98").Value = 1975.0416
$ws.Range("K98").Value = 1975.0416
$ws.Range("M98").Value = -477.0416
$ws.Range("H100").Value = 6442.757
$ws.Range("I100").Value = 1938.7693
$ws.Range("K100").Value = 1938.7693
$ws.Range("M100").Value = -1397.7693
$ws.Range("H103").Value = 2620
$ws.Range("J103").Value = 3800
$ws.Range("L103").Value = 11400
$ws.Range("N103").Value = -12572
$ws.Range("H111").Value = 15208.167
$ws.Range("I111").Value = 17956.223
$ws.Range("J111").Value = 6964
$ws.Range("K111").Value = 53868.66900000001
$ws.Range("L111").Value = 20892
$ws.Range("M111").Value = -50801.66900000001
$ws.Range("N111").Value = -27026
$ws.Range("H117").Value = 99595
$ws.Range("J117").Value = 99595
$ws.Range("L117").Value = 99595
$ws.Range("N117").Value = -108773
$ws.Range("H122").Value = 2284.8845
$ws.Range("I122").Value = 1975.0416
$ws.Range("K122").Value = 5925.1248
$ws.Range("M122").Value = -3475.1248
$ws.Range("H129").Value = 1635.289
$ws.Range("I129").Value = 713.6
$ws.Range("J129").Value = 1750.5
$ws.Range("K129").Value = 2140.8
$ws.Range("L129").Value = 5251.5
$ws.Range("M129").Value = 2859.2
$ws.Range("N129").Value = -15251.5
$ws.Range("H132").Value = 21280684
$ws.Range("I132").Value = 25644724
$ws.Range("K132").Value = 76934172
$ws.Range("M132").Value = -76931642
$ws.Range("H137").Value = 1805.5238
$ws.Range("I137").Value = 1299.3636
$ws.Range("J137").Value = 2362.3
$ws.Range("K137").Value = 3898.0908
$ws.Range("L137").Value = 7086.900000000001
$ws.Range("M137").Value = -1348.0908
$ws.Range("N137").Value = -12186.9
$ws.Range("H138").Value = 405624.8
$ws.Range("I138").Value = 2740.6365
$ws.Range("J138").Value = 722176.6
$ws.Range("K138").Value = 8221.9095
$ws.Range("L138").Value = 2166529.8
$ws.Range("M138").Value = -3081.9095
$ws.Range("N138").Value = -2176809.8
$ws.Range("H141").Value = 1760.25
$ws.Range("I141").Value = 1524.9714
$ws.Range("K141").Value = 4574.914199999999
$ws.Range("M141").Value = 605.0858000000007

# ---- Sheet: ARM (sheet2.xml) ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2163.0454
$ws.Range("I2").Value = 2167.282
$ws.Range("K2").Value = 2167.282
$ws.Range("M2").Value = -2054.282
$ws.Range("H5").Value = 170.76
$ws.Range("J5").Value = 63.666668
$ws.Range("L5").Value = 63.666668
$ws.Range("N5").Value = -287.666668
$ws.Range("H32").Value = 1652.5264
$ws.Range("I32").Value = 1573.25
$ws.Range("J32").Value = 2477
$ws.Range("K32").Value = 1573.25
$ws.Range("L32").Value = 2477
$ws.Range("M32").Value = -1286.25
$ws.Range("N32").Value = -3051
$ws.Range("H45").Value = 8224.5625
$ws.Range("I45").Value = 7867.5454
$ws.Range("K45").Value = 7867.5454
$ws.Range("M45").Value = -7490.5454
$ws.Range("H56").Value = 0
$ws.Range("J56").Value = 0
$ws.Range("L56").Value = 0
$ws.Range("N56").ClearContents()
$ws.Range("H61").Value = 22666
$ws.Range("I61").Value = 23199.2
$ws.Range("J61").Value = 20000
$ws.Range("K61").Value = 23199.2
$ws.Range("L61").Value = 20000
$ws.Range("M61").Value = -22987.2
$ws.Range("N61").Value = -20424
$ws.Range("H74").Value = 2393.8604
$ws.Range("I74").Value = 2004.3103
$ws.Range("K74").Value = 2004.3103
$ws.Range("M74").Value = -1130.3103
$ws.Range("H77").Value = 2393.8604
$ws.Range("I77").Value = 2004.3103
$ws.Range("K77").Value = 10021.5515
$ws.Range("M77").Value = -5653.551500000001
$ws.Range("H109").Value = 85000
$ws.Range("J109").Value = 85000
$ws.Range("L109").Value = 85000
$ws.Range("N109").Value = -87774
$ws.Range("H110").Value = 2281.0557
$ws.Range("I110").Value = 2139.818
$ws.Range("K110").Value = 2139.818
$ws.Range("M110").Value = -94.81800000000021
$ws.Range("H116").Value = 2163.0454
$ws.Range("I116").Value = 2167.282
$ws.Range("K116").Value = 2167.282
$ws.Range("M116").Value = 126.7179999999998
$ws.Range("H132").Value = 3621.818
$ws.Range("I132").Value = 3680.125
$ws.Range("J132").Value = 3466.3333
$ws.Range("K132").Value = 11040.375
$ws.Range("L132").Value = 10398.9999
$ws.Range("M132").Value = -8510.375
$ws.Range("N132").Value = -15458.9999
$ws.Range("H136").Value = 22666
$ws.Range("I136").Value = 23199.2
$ws.Range("J136").Value = 20000
$ws.Range("K136").Value = 69597.60000000001
$ws.Range("L136").Value = 60000
$ws.Range("M136").Value = -67047.60000000001
$ws.Range("N136").Value = -65100

# ---- Sheet: BSM (sheet3.xml) ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2163.0454
$ws.Range("I3").Value = 2167.282
$ws.Range("K3").Value = 2167.282
$ws.Range("M3").Value = -2053.282
$ws.Range("H4").Value = 170.76
$ws.Range("J4").Value = 63.666668
$ws.Range("L4").Value = 63.666668
$ws.Range("N4").Value = -293.666668
$ws.Range("H26").Value = 26435.5
$ws.Range("I26").Value = 26435.5
$ws.Range("K26").Value = 26435.5
$ws.Range("M26").Value = -26143.5
$ws.Range("H86").Value = 2679980.5
$ws.Range("I86").Value = 4782503.5
$ws.Range("K86").Value = 4782503.5
$ws.Range("M86").Value = -4781380.5
$ws.Range("H89").Value = 2679980.5
$ws.Range("I89").Value = 4782503.5
$ws.Range("K89").Value = 23912517.5
$ws.Range("M89").Value = -23906901.5
$ws.Range("H99").Value = 1173.8334
$ws.Range("I99").Value = 1213.8182
$ws.Range("K99").Value = 1213.8182
$ws.Range("M99").Value = 284.1818000000001
$ws.Range("H105").Value = 2962.973
$ws.Range("I105").Value = 1723.6818
$ws.Range("J105").Value = 4780.6
$ws.Range("K105").Value = 1723.6818
$ws.Range("L105").Value = 4780.6
$ws.Range("M105").Value = 23.31819999999993
$ws.Range("N105").Value = -8274.6
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("M115").ClearContents()
$ws.Range("N115").ClearContents()
$ws.Range("H134").Value = 3768.775
$ws.Range("I134").Value = 3300.5334
$ws.Range("K134").Value = 9901.600199999999
$ws.Range("M134").Value = -7366.600199999999

# ---- Sheet: CRP (sheet4.xml) ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 4002
$ws.Range("I2").Value = 4
$ws.Range("K2").Value = 4
$ws.Range("M2").Value = 109
$ws.Range("H31").Value = 3673.8845
$ws.Range("I31").Value = 2641.3333
$ws.Range("J31").Value = 4558.9287
$ws.Range("K31").Value = 2641.3333
$ws.Range("L31").Value = 4558.9287
$ws.Range("M31").Value = -2346.3333
$ws.Range("N31").Value = -5148.9287
$ws.Range("H34").Value = 3673.8845
$ws.Range("I34").Value = 2641.3333
$ws.Range("J34").Value = 4558.9287
$ws.Range("K34").Value = 2641.3333
$ws.Range("L34").Value = 4558.9287
$ws.Range("M34").Value = -2439.3333
$ws.Range("N34").Value = -4962.9287
$ws.Range("H43").Value = 43997.25
$ws.Range("J43").Value = 43997.25
$ws.Range("L43").Value = 43997.25
$ws.Range("N43").Value = -44365.25
$ws.Range("H58").Value = 2190.647
$ws.Range("I58").Value = 2682.875
$ws.Range("J58").Value = 1753.1111
$ws.Range("K58").Value = 2682.875
$ws.Range("L58").Value = 1753.1111
$ws.Range("M58").Value = -2479.875
$ws.Range("N58").Value = -2159.1111
$ws.Range("H94").Value = 1398.5625
$ws.Range("I94").Value = 1998
$ws.Range("J94").Value = 1312.9286
$ws.Range("K94").Value = 1998
$ws.Range("L94").Value = 1312.9286
$ws.Range("M94").Value = -1547
$ws.Range("N94").Value = -2214.9286
$ws.Range("H99").Value = 4413.2964
$ws.Range("I99").Value = 4302.909
$ws.Range("J99").Value = 4899
$ws.Range("K99").Value = 4302.909
$ws.Range("L99").Value = 4899
$ws.Range("M99").Value = -2804.909
$ws.Range("N99").Value = -7895
$ws.Range("H101").Value = 43997.25
$ws.Range("J101").Value = 43997.25
$ws.Range("L101").Value = 43997.25
$ws.Range("N101").Value = -50487.25
$ws.Range("H105").Value = 4881.1665
$ws.Range("I105").Value = 4123.5557
$ws.Range("J105").Value = 7154
$ws.Range("K105").Value = 4123.5557
$ws.Range("L105").Value = 7154
$ws.Range("M105").Value = -2376.5557
$ws.Range("N105").Value = -10648
$ws.Range("H122").Value = 3603.1052
$ws.Range("I122").Value = 1820.3846
$ws.Range("K122").Value = 5461.1538
$ws.Range("M122").Value = -3011.1538
$ws.Range("H126").Value = 4413.2964
$ws.Range("I126").Value = 4302.909
$ws.Range("J126").Value = 4899
$ws.Range("K126").Value = 12908.727
$ws.Range("L126").Value = 14697
$ws.Range("M126").Value = -10438.727
$ws.Range("N126").Value = -19637
$ws.Range("H132").Value = 1682.2
$ws.Range("I132").Value = 1852.75
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 5558.25
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -3028.25
$ws.Range("N132").Value = -8060
$ws.Range("H134").Value = 974.6667
$ws.Range("I134").Value = 974.6667
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 2924.0001
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -389.0001000000002
$ws.Range("N134").ClearContents()
$ws.Range("H136").Value = 2190.647
$ws.Range("I136").Value = 2682.875
$ws.Range("J136").Value = 1753.1111
$ws.Range("K136").Value = 8048.625
$ws.Range("L136").Value = 5259.3333
$ws.Range("M136").Value = -5498.625
$ws.Range("N136").Value = -10359.3333

# ---- Sheet: CUL (sheet5.xml) ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 36545450
$ws.Range("J4").Value = 8840748
$ws.Range("L4").Value = 26522244
$ws.Range("N4").Value = -26522468
$ws.Range("H41").Value = 5000
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("H60").Value = 1264.4
$ws.Range("I60").Value = 190.66667
$ws.Range("K60").Value = 572.00001
$ws.Range("M60").Value = -321.00001
$ws.Range("H62").Value = 5137.3335
$ws.Range("J62").Value = 5000
$ws.Range("L62").Value = 15000
$ws.Range("N62").Value = -16372
$ws.Range("H65").Value = 5137.3335
$ws.Range("J65").Value = 5000
$ws.Range("L65").Value = 45000
$ws.Range("N65").Value = -51864
$ws.Range("H68").Value = 2487.138
$ws.Range("I68").Value = 1376.1111
$ws.Range("J68").Value = 2691.204
$ws.Range("K68").Value = 4128.3333
$ws.Range("L68").Value = 8073.612000000001
$ws.Range("M68").Value = -3317.3333
$ws.Range("N68").Value = -9695.612000000001
$ws.Range("H71").Value = 2487.138
$ws.Range("I71").Value = 1376.1111
$ws.Range("J71").Value = 2691.204
$ws.Range("K71").Value = 12384.9999
$ws.Range("L71").Value = 24220.836
$ws.Range("M71").Value = -8328.999900000001
$ws.Range("N71").Value = -32332.836
$ws.Range("H98").Value = 874.8
$ws.Range("J98").Value = 918.5
$ws.Range("L98").Value = 2755.5
$ws.Range("N98").Value = -5751.5
$ws.Range("H107").Value = 2160.7646
$ws.Range("J107").Value = 2556.3845
$ws.Range("L107").Value = 7669.1535
$ws.Range("N107").Value = -11509.1535
$ws.Range("H112").Value = 254957.25
$ws.Range("J112").Value = 6610
$ws.Range("L112").Value = 19830
$ws.Range("N112").Value = -22046
$ws.Range("H122").Value = 821.43475
$ws.Range("I122").Value = 467.625
$ws.Range("J122").Value = 1010.13336
$ws.Range("K122").Value = 4208.625
$ws.Range("L122").Value = 9091.20024
$ws.Range("M122").Value = -1758.625
$ws.Range("N122").Value = -13991.20024
$ws.Range("H131").Value = 1897171.4
$ws.Range("J131").Value = 2274450
$ws.Range("L131").Value = 6823350
$ws.Range("N131").Value = -6833430
$ws.Range("H140").Value = 2585.575
$ws.Range("I140").Value = 2298.808
$ws.Range("K140").Value = 6896.424
$ws.Range("M140").Value = -1716.424
$ws.Range("H141").Value = 28606
$ws.Range("I141").Value = 28606
$ws.Range("K141").Value = 85818
$ws.Range("M141").Value = -80638

# ---- Sheet: GSM (sheet6.xml) ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 14998
$ws.Range("J59").Value = 14998
$ws.Range("L59").Value = 14998
$ws.Range("N59").Value = -16164
$ws.Range("H97").Value = 1382.4706
$ws.Range("I97").Value = 1285.2142
$ws.Range("K97").Value = 1285.2142
$ws.Range("M97").Value = -789.2141999999999
$ws.Range("N97").Value = -2828.3334
$ws.Range("H102").Value = 27193.6
$ws.Range("I102").Value = 1161.3572
$ws.Range("K102").Value = 1161.3572
$ws.Range("M102").Value = 460.6428000000001
$ws.Range("H122").Value = 57250.156
$ws.Range("I122").Value = 75498.36
$ws.Range("J122").Value = 6155.2
$ws.Range("K122").Value = 226495.08
$ws.Range("L122").Value = 18465.6
$ws.Range("M122").Value = -224045.08
$ws.Range("N122").Value = -23365.6
$ws.Range("H123").Value = 42206.066
$ws.Range("J123").Value = 42206.066
$ws.Range("L123").Value = 42206.066
$ws.Range("N123").Value = -47106.066
$ws.Range("H126").Value = 70371.60000000001
$ws.Range("I126").Value = 79658.53999999999
$ws.Range("K126").Value = 238975.62
$ws.Range("M126").Value = -236505.62
$ws.Range("H132").Value = 7577.4116
$ws.Range("I132").Value = 9196.308000000001
$ws.Range("J132").Value = 2316
$ws.Range("K132").Value = 27588.924
$ws.Range("L132").Value = 6948
$ws.Range("M132").Value = -25058.924
$ws.Range("N132").Value = -12008

# ---- Sheet: LTW (sheet7.xml) ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 22698.4
$ws.Range("I7").Value = 25314.875
$ws.Range("J7").Value = 19708.143
$ws.Range("K7").Value = 25314.875
$ws.Range("L7").Value = 19708.143
$ws.Range("M7").Value = -25202.875
$ws.Range("N7").Value = -19932.143
$ws.Range("H16").Value = 2173.6775
$ws.Range("I16").Value = 1073.5333
$ws.Range("J16").Value = 3205.0625
$ws.Range("K16").Value = 1073.5333
$ws.Range("L16").Value = 3205.0625
$ws.Range("M16").Value = -903.5333000000001
$ws.Range("N16").Value = -3545.0625
$ws.Range("H40").Value = 6867.8823
$ws.Range("I40").Value = 7066.5386
$ws.Range("K40").Value = 7066.5386
$ws.Range("M40").Value = -6930.5386
$ws.Range("H61").Value = 3111.5518
$ws.Range("J61").Value = 6398
$ws.Range("L61").Value = 6398
$ws.Range("N61").Value = -6802
$ws.Range("H64").Value = 37599.168
$ws.Range("J64").Value = 37599.168
$ws.Range("L64").Value = 37599.168
$ws.Range("N64").Value = -38049.168
$ws.Range("H67").Value = 37599.168
$ws.Range("J67").Value = 37599.168
$ws.Range("L67").Value = 37599.168
$ws.Range("N67").Value = -39159.168
$ws.Range("H68").Value = 4334.6665
$ws.Range("I68").Value = 5499.5
$ws.Range("J68").Value = 3402.8
$ws.Range("K68").Value = 5499.5
$ws.Range("L68").Value = 3402.8
$ws.Range("M68").Value = -4750.5
$ws.Range("N68").Value = -4900.8
$ws.Range("H71").Value = 4334.6665
$ws.Range("I71").Value = 5499.5
$ws.Range("J71").Value = 3402.8
$ws.Range("K71").Value = 27497.5
$ws.Range("L71").Value = 17014
$ws.Range("M71").Value = -23753.5
$ws.Range("N71").Value = -24502
$ws.Range("H74").Value = 26764.705
$ws.Range("J74").Value = 25000
$ws.Range("L74").Value = 25000
$ws.Range("N74").Value = -26996
$ws.Range("H77").Value = 26764.705
$ws.Range("J77").Value = 25000
$ws.Range("L77").Value = 75000
$ws.Range("N77").Value = -84984
$ws.Range("H82").Value = 4959.8887
$ws.Range("I82").Value = 7684
$ws.Range("J82").Value = 1554.75
$ws.Range("K82").Value = 7684
$ws.Range("L82").Value = 1554.75
$ws.Range("M82").Value = -7323
$ws.Range("N82").Value = -2276.75
$ws.Range("H85").Value = 4959.8887
$ws.Range("I85").Value = 7684
$ws.Range("J85").Value = 1554.75
$ws.Range("K85").Value = 7684
$ws.Range("L85").Value = 1554.75
$ws.Range("M85").Value = -6436
$ws.Range("N85").Value = -4050.75
$ws.Range("H87").Value = 89990
$ws.Range("J87").Value = 89990
$ws.Range("L87").Value = 89990
$ws.Range("N87").Value = -92236
$ws.Range("H90").Value = 89990
$ws.Range("J90").Value = 89990
$ws.Range("L90").Value = 269970
$ws.Range("N90").Value = -281202
$ws.Range("H93").Value = 1273.8889
$ws.Range("I93").Value = 1066.4286
$ws.Range("K93").Value = 1066.4286
$ws.Range("M93").Value = 181.5714
$ws.Range("H102").Value = 96420
$ws.Range("J102").Value = 96420
$ws.Range("L102").Value = 96420
$ws.Range("N102").Value = -102910
$ws.Range("H113").Value = 3111.5518
$ws.Range("J113").Value = 6398
$ws.Range("L113").Value = 6398
$ws.Range("N113").Value = -10738
$ws.Range("H122").Value = 17250.5
$ws.Range("I122").Value = 15251.25
$ws.Range("K122").Value = 45753.75
$ws.Range("M122").Value = -43303.75
$ws.Range("H126").Value = 22698.4
$ws.Range("I126").Value = 25314.875
$ws.Range("J126").Value = 19708.143
$ws.Range("K126").Value = 75944.625
$ws.Range("L126").Value = 59124.429
$ws.Range("M126").Value = -73474.625
$ws.Range("N126").Value = -64064.429
$ws.Range("H132").Value = 5130.6562
$ws.Range("I132").Value = 5274.625
$ws.Range("J132").Value = 4698.75
$ws.Range("K132").Value = 15823.875
$ws.Range("L132").Value = 14096.25
$ws.Range("M132").Value = -13293.875
$ws.Range("N132").Value = -19156.25
$ws.Range("H136").Value = 4135.763
$ws.Range("I136").Value = 4122.3706
$ws.Range("J136").Value = 4168.636
$ws.Range("K136").Value = 12367.1118
$ws.Range("L136").Value = 12505.908
$ws.Range("M136").Value = -9817.111800000001
$ws.Range("N136").Value = -17605.908
$ws.Range("H137").Value = 67096.875
$ws.Range("J137").Value = 67096.875
$ws.Range("L137").Value = 67096.875
$ws.Range("N137").Value = -77296.875

# ---- Sheet: WVR (sheet8.xml) ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 22844.715
$ws.Range("J45").Value = 24943.6
$ws.Range("L45").Value = 24943.6
$ws.Range("N45").Value = -25925.6
$ws.Range("H62").Value = 6500
$ws.Range("I62").Value = 6333.3335
$ws.Range("K62").Value = 6333.3335
$ws.Range("M62").Value = -5709.3335
$ws.Range("H65").Value = 6500
$ws.Range("I65").Value = 6333.3335
$ws.Range("K65").Value = 31666.6675
$ws.Range("M65").Value = -28546.6675
$ws.Range("H68").Value = 49999
$ws.Range("J68").Value = 49999
$ws.Range("L68").Value = 49999
$ws.Range("N68").Value = -51621
$ws.Range("H71").Value = 49999
$ws.Range("J71").Value = 49999
$ws.Range("L71").Value = 149997
$ws.Range("N71").Value = -158109
$ws.Range("H81").Value = 85176.30499999999
$ws.Range("I81").Value = 205238.8
$ws.Range("J81").Value = 10137.25
$ws.Range("K81").Value = 410477.6
$ws.Range("L81").Value = 20274.5
$ws.Range("M81").Value = -409416.6
$ws.Range("N81").Value = -22396.5
$ws.Range("H84").Value = 85176.30499999999
$ws.Range("I84").Value = 205238.8
$ws.Range("J84").Value = 10137.25
$ws.Range("K84").Value = 2052388
$ws.Range("L84").Value = 101372.5
$ws.Range("M84").Value = -2047084
$ws.Range("N84").Value = -111980.5
$ws.Range("H122").Value = 3211.8928
$ws.Range("I122").Value = 951.55
$ws.Range("K122").Value = 2854.65
$ws.Range("M122").Value = -404.6499999999996
$ws.Range("H132").Value = 1225.7667
$ws.Range("I132").Value = 991.8929000000001
$ws.Range("J132").Value = 4500
$ws.Range("K132").Value = 2975.6787
$ws.Range("L132").Value = 13500
$ws.Range("M132").Value = -445.6787000000004
$ws.Range("N132").Value = -18560
$ws.Range("H136").Value = 7788.2915
$ws.Range("I136").Value = 9584.134
$ws.Range("K136").Value = 28752.402
$ws.Range("M136").Value = -26202.402
